$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Row 4 (DigitalDownloadsExecution test data row) is being updated:
#  - TestCase name changes from "...checkout" to "...digitalDownloads"
#  - Email used for the address form changes to "meghadmegha@gmail.com"
#  - New assertion columns (O-X) are populated to mirror the entered
#    address details and confirm the login assertion text.

$ws.Range("B4").Value = "meghadevaraja1998@gmail.com"
$ws.Range("C4").Value = "Megha@123"
$ws.Range("F4").Value = "Megha"
$ws.Range("G4").Value = "Devaraj"
$ws.Range("H4").Value = "meghadmegha@gmail.com"
$ws.Range("I4").Value = "Bengaluru"
$ws.Range("J4").Value = "kattreguppe"
$ws.Range("K4").Value = 560085
$ws.Range("L4").Value = 9874563321
$ws.Range("M4").Value = "India"
$ws.Range("N4").Value = "New Address"
$ws.Range("A4").Value = "DigitalDownloadsExecution.digitalDownloads"
$ws.Range("O4").Value = "Megha"
$ws.Range("P4").Value = "Devaraj"
$ws.Range("Q4").Value = "meghadmegha@gmail.com"
$ws.Range("R4").Value = "New Address"
$ws.Range("S4").Value = "India"
$ws.Range("T4").Value = "Bengaluru"
$ws.Range("U4").Value = "kattreguppe"
$ws.Range("V4").Value = 560085
$ws.Range("W4").Value = 9874563321
$ws.Range("X4").Value = "You logged into a secure area!"

# The old row 4 had mailto hyperlinks on B4, C4 and H4 (pointing at the
# former username/password/email). These no longer apply and are removed.
$ws.Range("B4").Hyperlinks.Delete()
$ws.Range("C4").Hyperlinks.Delete()
$ws.Range("H4").Hyperlinks.Delete()
